$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the sheet is being substantially restructured/expanded.
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# Header row (row 1): one label per column, A..S
#
# Values are written in the same order the shared-string table ends up in
# (the original 12 strings first, in their original order, then the 7 new
# ones in the order: ErrNV, ErrNum, Err3Num, Err3NV, Err2Num, Err2NV,
# ErrMixed) so the rebuilt sharedStrings.xml matches the source workbook.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "ErrDiv0"
$ws.Range("B1").Value = "ErrName"
$ws.Range("C1").Value = "ErrValue"
$ws.Range("D1").Value = "ErrRef"
$ws.Range("G1").Value = "Err2Div0"
$ws.Range("H1").Value = "Err2Name"
$ws.Range("I1").Value = "Err2Value"
$ws.Range("J1").Value = "Err2Ref"
$ws.Range("M1").Value = "Err3Div0"
$ws.Range("N1").Value = "Err3Name"
$ws.Range("O1").Value = "Err3Value"
$ws.Range("P1").Value = "Err3Ref"
$ws.Range("E1").Value = "ErrNV"
$ws.Range("F1").Value = "ErrNum"
$ws.Range("R1").Value = "Err3Num"
$ws.Range("Q1").Value = "Err3NV"
$ws.Range("L1").Value = "Err2Num"
$ws.Range("K1").Value = "Err2NV"
$ws.Range("S1").Value = "ErrMixed"

# ---------------------------------------------------------------------------
# Group 1 (columns A-F): formula lives in row 2, row 3 holds a plain 1
# ---------------------------------------------------------------------------
$ws.Range("A2").Formula = "=8/0"
$ws.Range("B2").Formula = "=a0"
$ws.Range("C2").Formula = "=INT(""kjk"")"
$ws.Range("D2").Formula = "=#REF!"
$ws.Range("E2").Formula = "=VLOOKUP(""1"",D1:E1,1)"

# F2:F3 are written as a shared formula pair, then F3 is overwritten with a
# plain value (mirrors how the workbook stores it).
$ws.Range("F2:F3").Formula = "=SQRT(-4)"
$ws.Range("F3").Value = 1

$ws.Range("A3:E3").Value = 1

# ---------------------------------------------------------------------------
# Group 2 (columns G-L): mirror of group 1 but formula lives in row 3
# ---------------------------------------------------------------------------
$ws.Range("G2:K2").Value = 1

$ws.Range("G3").Formula = "=8/0"
$ws.Range("H3").Formula = "=a0"
$ws.Range("I3").Formula = "=INT(""kjk"")"
$ws.Range("J3").Formula = "=#REF!"
$ws.Range("K3").Formula = "=VLOOKUP(""1"",J2:K2,1)"

# ---------------------------------------------------------------------------
# Group 3 (columns M-R): formula repeated on both row 2 and row 3
# ---------------------------------------------------------------------------
$ws.Range("M2:M3").Formula = "=8/0"
$ws.Range("N2:N3").Formula = "=a0"
$ws.Range("O2:O3").Formula = "=INT(""kjk"")"

$ws.Range("P2").Formula = "=#REF!"
$ws.Range("P3").Formula = "=#REF!"

$ws.Range("Q2").Formula = "=VLOOKUP(""1"",P1:Q1,1)"
$ws.Range("Q3").Formula = "=VLOOKUP(""1"",P2:Q2,1)"

$ws.Range("R2:R3").Formula = "=SQRT(-4)"

# L2 is a plain value; L3 holds the formula (mirrors columns G-K in this
# group, where the formula lives in row 3).
$ws.Range("L2").Value = 1
$ws.Range("L3").Formula = "=SQRT(-4)"

# ---------------------------------------------------------------------------
# Column S ("ErrMixed"): one error type per row, rows 3-8 (row 2 is a plain 1)
# ---------------------------------------------------------------------------
$ws.Range("S2").Value = 1
$ws.Range("S3").Formula = "=8/0"
$ws.Range("S4").Formula = "=a0"
$ws.Range("S5").Formula = "=INT(""kjk"")"
$ws.Range("S6").Formula = "=#REF!"
$ws.Range("S7").Formula = "=VLOOKUP(""1"",R6:S6,1)"
$ws.Range("S8").Formula = "=SQRT(-4)"

# ---------------------------------------------------------------------------
# Selection state, matching the saved workbook
# ---------------------------------------------------------------------------
$ws.Range("S4").Select()
